$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: rename / reorder the header labels -----------------------------
# Old:  A StaffCreationId | B LeaveTypeName | C TransactionFlag | D LeaveCount
#       E LeaveReason     | F Month         | G Year            | H Remarks
# New:  A StaffId | B LeaveType | C TransactionFlag | D LeaveReason
#       E Month    | F Year      | G LeaveCount      | H Remarks
$ws.Range("A1").Value = "StaffId"
$ws.Range("B1").Value = "LeaveType"
$ws.Range("C1").Value = "TransactionFlag"
$ws.Range("D1").Value = "LeaveReason"
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Year"
$ws.Range("G1").Value = "LeaveCount"
$ws.Range("H1").Value = "Remarks"

# --- Row 2: drop the old sample/template data row ---------------------------
$ws.Range("A2:H2").ClearContents()

# --- Leftover formatted cells further out on row 2 --------------------------
$ws.Range("M2").NumberFormat = "mm:ss.0"
$ws.Range("O2").NumberFormat = "mm:ss.0"

# --- Move the active selection ----------------------------------------------
$ws.Range("H4").Select() | Out-Null
